$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" header text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 01:35"

# --- Simple statistic refreshes (country unchanged, only numbers move) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 1706150
$ws.Range("C4").Value = 19714
$ws.Range("D4").Value = 462757
$ws.Range("E4").Value = 1143589
$ws.Range("G4").Value = 504
$ws.Range("H4").Value = 99804

# Row 13: India
$ws.Range("B13").Value = 144950
$ws.Range("C13").Value = 6414
$ws.Range("E13").Value = 80072

# Row 47: Argentina
$ws.Range("B47").Value = 12628
$ws.Range("C47").Value = 552
$ws.Range("E47").Value = 8162
$ws.Range("G47").Value = 15
$ws.Range("H47").Value = 467

# Row 137: Reunion
$ws.Range("B137").Value = 456
$ws.Range("C137").Value = 4
$ws.Range("E137").Value = 44

# --- Re-sort "Guayana Francesa" into its alphabetical slot right after "Togo" ---
# It previously sat right before "Vietnam" (after "Mauricio"). Effectively this
# removes it from its old spot and inserts a fresh row for it just after "Togo",
# shifting "Ruanda", "Isla de Man" and "Mauricio" down by one row with their
# pre-existing stats, while "Guayana Francesa" gets freshly updated numbers.

# Row 142 becomes Guayana Francesa with new stats
$ws.Range("A142").Value = "Guayana Francesa"
$ws.Range("B142").Value = 353
$ws.Range("C142").Value = 25
$ws.Range("D142").Value = 146
$ws.Range("E142").Value = 206
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 1

# Row 143 becomes Ruanda (old row 142 stats)
$ws.Range("A143").Value = "Ruanda"
$ws.Range("B143").Value = 336
$ws.Range("C143").Value = 9
$ws.Range("D143").Value = 238
$ws.Range("E143").Value = 98
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

# Row 144 becomes Isla de Man (old row 143 stats)
$ws.Range("A144").Value = "Isla de Man"
$ws.Range("B144").Value = 336
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 305
$ws.Range("E144").Value = 7
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 24

# Row 145 becomes Mauricio (old row 144 stats)
$ws.Range("A145").Value = "Mauricio"
$ws.Range("B145").Value = 334
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 322
$ws.Range("E145").Value = 2
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 10
